$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump
$meta.Range("B3").Value = "6.0.0"

# Date update
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 "Contact" -> "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Remove the duplicate "Contact" row (old row 11), shifting remaining rows up
$meta.Rows.Item(11).Delete()

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the root "Extension" element - give it the real Short/Definition text
$elements.Range("K2").Value = "Employee Pay Frequency"
$elements.Range("L2").Value = "Code indicating the frequency with which the employee is paid (e.g., daily, weekly, monthly)"
